# Continue digitizing: add "morici" study rows (fwd / cwd components) to the
# "data" worksheet, and fix the thin_type ("J") column for existing fwd/cwd
# rows that should be "NA" (this variable doesn't apply to those records).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1) Fix thin_type (column J) on pre-existing rows: these rows record fwd/cwd
#    data where thin_type doesn't apply, so it should read "NA" instead of
#    "shelterwood"/"commercial".
# ---------------------------------------------------------------------------
$naRows = @(14,15,26,27,38,39,50,51,62,63,74,75,76,86,87,88,98,99,100,110,111,112,122,123,124,134,135,136,146,147,148,158,159,160,170,171,178,179)
foreach ($r in $naRows) {
    $ws.Cells.Item($r, 10).Value = "NA"
}

# ---------------------------------------------------------------------------
# 2) Append new rows 186-209 with the newly-digitized "morici" study data.
#    Columns: A value, B variable, C units, D years_post, E treatment,
#    F study, G forest_type, H region, I burn_season, J thin_type, K notes.
#    Values are written column-by-column (matching how the workbook was
#    actually authored) so new shared strings land in the same order.
# ---------------------------------------------------------------------------
$startRow = 186
$nRows = 24

$valueCol = @(23.9, 19.3, 25, 17.100000000000001, 15.3, 19.899999999999999, `
              10, 8.4, 3.7, 7.4, 6.9, 16.399999999999999, `
              6.7, 3.9, 4.8, 4.4000000000000004, 3.1, 3.2, `
              1.96, 0.76, 1.31, 0.27, 0.05, 0.27)

$variableCol = @("all_woody","all_woody","all_woody", `
                 "shelterwood","shelterwood","shelterwood", `
                 "cwd_sound","cwd_sound","cwd_sound", `
                 "cwd_rotten","cwd_rotten","cwd_rotten", `
                 "fwd","fwd","fwd", `
                 "hundred_hour","hundred_hour","hundred_hour", `
                 "ten_hour","ten_hour","ten_hour", `
                 "one_hour","one_hour","one_hour")

$yearsPostCol = @(0,8,17,0,8,17,0,8,17,0,8,17,0,8,17,0,8,17,0,8,17,0,8,17)

for ($i = 0; $i -lt $nRows; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $valueCol[$i]
}
for ($i = 0; $i -lt $nRows; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $variableCol[$i]
}
for ($i = 0; $i -lt $nRows; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = "Mg/ha"
}
for ($i = 0; $i -lt $nRows; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $yearsPostCol[$i]
}
for ($i = 0; $i -lt $nRows; $i++) {
    $ws.Cells.Item($startRow + $i, 5).Value = "control"
}
for ($i = 0; $i -lt $nRows; $i++) {
    $ws.Cells.Item($startRow + $i, 6).Value = "morici"
}
for ($i = 0; $i -lt $nRows; $i++) {
    $ws.Cells.Item($startRow + $i, 7).Value = "ponderosa"
}
for ($i = 0; $i -lt $nRows; $i++) {
    $ws.Cells.Item($startRow + $i, 8).Value = "interior_pnw"
}
for ($i = 0; $i -lt $nRows; $i++) {
    $ws.Cells.Item($startRow + $i, 9).Value = "NA"
}
for ($i = 0; $i -lt $nRows; $i++) {
    $ws.Cells.Item($startRow + $i, 10).Value = "NA"
}
for ($i = 0; $i -lt $nRows; $i++) {
    $ws.Cells.Item($startRow + $i, 11).Value = "NA"
}

# ---------------------------------------------------------------------------
# 3) Update the view so the newly-added rows are visible / selected.
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 171
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("G180").Select() | Out-Null
